$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- New header/value pairs in F1:G2, mirroring the existing N7:O7 (header)
#     and N8:O8 (data) cells for Battery Standby / Alarm Load detail columns ---

# F1/G1: header style + text (copy format from N7/O7, then set the text)
$ws.Range("N7").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "AlarmLoadingDetail"

$ws.Range("O7").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "StandbyLoadingDetail"

# F2/G2: data style + text (copy format from N8/O8, then set the text)
$ws.Range("N8").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = "Battery Alarm (A)"

$ws.Range("O8").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "Battery Standby (A)"

# Widen column G (bestFit-style) to accommodate the new "StandbyLoadingDetail" header
$ws.Columns.Item(7).ColumnWidth = 18.91

# Update the view: select F1:G2 (also clears the stale topLeftCell="K1" scroll anchor)
$ws.Range("F1:G2").Select() | Out-Null
